$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IoT_color_detection_BOM")

# The BOM used separate SMD_1207..SMD_1212 footprint names for the R30:R35
# resistor rows; replace them all with the common SMD_1206 footprint
# (already used elsewhere in the sheet, e.g. D29).
$ws.Range("D30:D35").Value = "SMD_1206"

# Re-assigning .Value resets quote-prefix formatting on these cells; restore
# the original column-D look by copying formats from an untouched D cell
# that already carries the SMD_1206 value/style.
$ws.Range("D29").Copy()
$ws.Range("D30:D35").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view: zoom to 130%, selection on D33 (matches the saved
# workbook state after the edit).
$ws.Select()
$excel.ActiveWindow.Zoom = 130
$ws.Range("D33").Select()
